# Applies "added harvard case classification" edit:
#  - Swaps the BP1/BQ1 header labels ("average_doctor" <-> "average_doctor_old")
#  - Updates statistic values across the sheet (rows 4-13) to reflect the
#    reclassified data, including shifting the former BP (average_doctor)
#    values into BQ (average_doctor_old) and writing freshly computed
#    averages into BP.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label swap (row 1)
$ws.Cells.Item(1, 68).Value = "average_doctor_old"  # BP1
$ws.Cells.Item(1, 69).Value = "average_doctor"  # BQ1

# Updated statistic values (rows 4-13)
# Row 4
$ws.Cells.Item(4, 5).Value = 0.478  # E4
$ws.Cells.Item(4, 6).Value = 0.052  # F4
$ws.Cells.Item(4, 7).Value = 0.229  # G4
$ws.Cells.Item(4, 14).Value = 0.474  # N4
$ws.Cells.Item(4, 15).Value = 0.061  # O4
$ws.Cells.Item(4, 16).Value = 0.248  # P4
$ws.Cells.Item(4, 17).Value = 0.051  # Q4
$ws.Cells.Item(4, 18).Value = 0.034  # R4
$ws.Cells.Item(4, 19).Value = 0.186  # S4
$ws.Cells.Item(4, 23).Value = 0.366  # W4
$ws.Cells.Item(4, 35).Value = 0.401  # AI4
$ws.Cells.Item(4, 36).Value = 0.093  # AJ4
$ws.Cells.Item(4, 37).Value = 0.306  # AK4
$ws.Cells.Item(4, 47).Value = 0.239  # AU4
$ws.Cells.Item(4, 49).Value = 0.161  # AW4
$ws.Cells.Item(4, 53).Value = 2.064  # BA4
$ws.Cells.Item(4, 54).Value = 0.143  # BB4
$ws.Cells.Item(4, 55).Value = 0.379  # BC4
$ws.Cells.Item(4, 59).Value = 0.729  # BG4
$ws.Cells.Item(4, 60).Value = 0.141  # BH4
$ws.Cells.Item(4, 61).Value = 0.375  # BI4
$ws.Cells.Item(4, 65).Value = 0.756  # BM4
$ws.Cells.Item(4, 66).Value = 0.064  # BN4
$ws.Cells.Item(4, 67).Value = 0.253  # BO4
$ws.Cells.Item(4, 68).Value = 0.688  # BP4
$ws.Cells.Item(4, 69).Value = 0.764  # BQ4
# Row 5
$ws.Cells.Item(5, 5).Value = 0.599  # E5
$ws.Cells.Item(5, 6).Value = 0.059  # F5
$ws.Cells.Item(5, 7).Value = 0.242  # G5
$ws.Cells.Item(5, 14).Value = 0.716  # N5
$ws.Cells.Item(5, 15).Value = 0.076  # O5
$ws.Cells.Item(5, 16).Value = 0.275  # P5
$ws.Cells.Item(5, 17).Value = 0.034  # Q5
$ws.Cells.Item(5, 18).Value = 0.014  # R5
$ws.Cells.Item(5, 19).Value = 0.119  # S5
$ws.Cells.Item(5, 23).Value = 0.335  # W5
$ws.Cells.Item(5, 24).Value = 0.097  # X5
$ws.Cells.Item(5, 25).Value = 0.312  # Y5
$ws.Cells.Item(5, 35).Value = 0.402  # AI5
$ws.Cells.Item(5, 36).Value = 0.09  # AJ5
$ws.Cells.Item(5, 37).Value = 0.3  # AK5
$ws.Cells.Item(5, 47).Value = 0.445  # AU5
$ws.Cells.Item(5, 48).Value = 0.077  # AV5
$ws.Cells.Item(5, 49).Value = 0.278  # AW5
$ws.Cells.Item(5, 53).Value = 1.3  # BA5
$ws.Cells.Item(5, 54).Value = 0.072  # BB5
$ws.Cells.Item(5, 55).Value = 0.268  # BC5
$ws.Cells.Item(5, 59).Value = 0.382  # BG5
$ws.Cells.Item(5, 60).Value = 0.05  # BH5
$ws.Cells.Item(5, 61).Value = 0.223  # BI5
$ws.Cells.Item(5, 65).Value = 0.521  # BM5
$ws.Cells.Item(5, 66).Value = 0.047  # BN5
$ws.Cells.Item(5, 67).Value = 0.216  # BO5
$ws.Cells.Item(5, 68).Value = 0.433  # BP5
$ws.Cells.Item(5, 69).Value = 0.458  # BQ5
# Row 6
$ws.Cells.Item(6, 5).Value = 0.532  # E6
$ws.Cells.Item(6, 14).Value = 0.57  # N6
$ws.Cells.Item(6, 17).Value = 0.041  # Q6
$ws.Cells.Item(6, 23).Value = 0.35  # W6
$ws.Cells.Item(6, 35).Value = 0.401  # AI6
$ws.Cells.Item(6, 47).Value = 0.311  # AU6
$ws.Cells.Item(6, 53).Value = 1.589  # BA6
$ws.Cells.Item(6, 59).Value = 0.501  # BG6
$ws.Cells.Item(6, 65).Value = 0.617  # BM6
$ws.Cells.Item(6, 68).Value = 0.53  # BP6
$ws.Cells.Item(6, 69).Value = 0.57  # BQ6
# Row 7
$ws.Cells.Item(7, 5).Value = 0.57  # E7
$ws.Cells.Item(7, 14).Value = 0.65  # N7
$ws.Cells.Item(7, 17).Value = 0.036  # Q7
$ws.Cells.Item(7, 23).Value = 0.341  # W7
$ws.Cells.Item(7, 35).Value = 0.402  # AI7
$ws.Cells.Item(7, 47).Value = 0.38  # AU7
$ws.Cells.Item(7, 53).Value = 1.402  # BA7
$ws.Cells.Item(7, 59).Value = 0.422  # BG7
$ws.Cells.Item(7, 65).Value = 0.556  # BM7
$ws.Cells.Item(7, 68).Value = 0.467  # BP7
$ws.Cells.Item(7, 69).Value = 0.497  # BQ7
# Row 8
$ws.Cells.Item(8, 5).Value = 0.706  # E8
$ws.Cells.Item(8, 6).Value = 0.07  # F8
$ws.Cells.Item(8, 7).Value = 0.265  # G8
$ws.Cells.Item(8, 14).Value = 0.802  # N8
$ws.Cells.Item(8, 15).Value = 0.058  # O8
$ws.Cells.Item(8, 16).Value = 0.242  # P8
$ws.Cells.Item(8, 17).Value = 0.037  # Q8
$ws.Cells.Item(8, 23).Value = 0.405  # W8
$ws.Cells.Item(8, 24).Value = 0.119  # X8
$ws.Cells.Item(8, 25).Value = 0.345  # Y8
$ws.Cells.Item(8, 35).Value = 0.472  # AI8
$ws.Cells.Item(8, 36).Value = 0.137  # AJ8
$ws.Cells.Item(8, 37).Value = 0.37  # AK8
$ws.Cells.Item(8, 47).Value = 0.394  # AU8
$ws.Cells.Item(8, 48).Value = 0.083  # AV8
$ws.Cells.Item(8, 49).Value = 0.288  # AW8
$ws.Cells.Item(8, 53).Value = 1.771  # BA8
$ws.Cells.Item(8, 54).Value = 0.108  # BB8
$ws.Cells.Item(8, 55).Value = 0.328  # BC8
$ws.Cells.Item(8, 59).Value = 0.564  # BG8
$ws.Cells.Item(8, 60).Value = 0.108  # BH8
$ws.Cells.Item(8, 61).Value = 0.329  # BI8
$ws.Cells.Item(8, 65).Value = 0.675  # BM8
$ws.Cells.Item(8, 66).Value = 0.061  # BN8
$ws.Cells.Item(8, 67).Value = 0.247  # BO8
$ws.Cells.Item(8, 68).Value = 0.59  # BP8
$ws.Cells.Item(8, 69).Value = 0.625  # BQ8
# Row 9
$ws.Cells.Item(9, 5).Value = 0.674  # E9
$ws.Cells.Item(9, 6).Value = 0.22  # F9
$ws.Cells.Item(9, 7).Value = 0.469  # G9
$ws.Cells.Item(9, 14).Value = 0.744  # N9
$ws.Cells.Item(9, 15).Value = 0.19  # O9
$ws.Cells.Item(9, 16).Value = 0.436  # P9
$ws.Cells.Item(9, 23).Value = 0.302  # W9
$ws.Cells.Item(9, 24).Value = 0.211  # X9
$ws.Cells.Item(9, 25).Value = 0.459  # Y9
$ws.Cells.Item(9, 35).Value = 0.419  # AI9
$ws.Cells.Item(9, 36).Value = 0.243  # AJ9
$ws.Cells.Item(9, 37).Value = 0.493  # AK9
$ws.Cells.Item(9, 53).Value = 1.744  # BA9
$ws.Cells.Item(9, 54).Value = 0.25  # BB9
$ws.Cells.Item(9, 55).Value = 0.5  # BC9
$ws.Cells.Item(9, 59).Value = 0.605  # BG9
$ws.Cells.Item(9, 60).Value = 0.239  # BH9
$ws.Cells.Item(9, 61).Value = 0.489  # BI9
$ws.Cells.Item(9, 65).Value = 0.651  # BM9
$ws.Cells.Item(9, 66).Value = 0.227  # BN9
$ws.Cells.Item(9, 67).Value = 0.477  # BO9
$ws.Cells.Item(9, 68).Value = 0.581  # BP9
$ws.Cells.Item(9, 69).Value = 0.618  # BQ9
# Row 10
$ws.Cells.Item(10, 5).Value = 0.814  # E10
$ws.Cells.Item(10, 6).Value = 0.151  # F10
$ws.Cells.Item(10, 7).Value = 0.389  # G10
$ws.Cells.Item(10, 14).Value = 0.93  # N10
$ws.Cells.Item(10, 15).Value = 0.065  # O10
$ws.Cells.Item(10, 16).Value = 0.255  # P10
$ws.Cells.Item(10, 23).Value = 0.512  # W10
$ws.Cells.Item(10, 24).Value = 0.25  # X10
$ws.Cells.Item(10, 25).Value = 0.5  # Y10
$ws.Cells.Item(10, 35).Value = 0.512  # AI10
$ws.Cells.Item(10, 36).Value = 0.25  # AJ10
$ws.Cells.Item(10, 37).Value = 0.5  # AK10
$ws.Cells.Item(10, 47).Value = 0.395  # AU10
$ws.Cells.Item(10, 48).Value = 0.239  # AV10
$ws.Cells.Item(10, 49).Value = 0.489  # AW10
$ws.Cells.Item(10, 53).Value = 2.186  # BA10
$ws.Cells.Item(10, 54).Value = 0.211  # BB10
$ws.Cells.Item(10, 55).Value = 0.459  # BC10
$ws.Cells.Item(10, 59).Value = 0.674  # BG10
$ws.Cells.Item(10, 60).Value = 0.22  # BH10
$ws.Cells.Item(10, 61).Value = 0.469  # BI10
$ws.Cells.Item(10, 65).Value = 0.814  # BM10
$ws.Cells.Item(10, 66).Value = 0.151  # BN10
$ws.Cells.Item(10, 67).Value = 0.389  # BO10
$ws.Cells.Item(10, 68).Value = 0.729  # BP10
$ws.Cells.Item(10, 69).Value = 0.758  # BQ10
# Row 11
$ws.Cells.Item(11, 5).Value = 0.86  # E11
$ws.Cells.Item(11, 6).Value = 0.12  # F11
$ws.Cells.Item(11, 7).Value = 0.347  # G11
$ws.Cells.Item(11, 14).Value = 0.93  # N11
$ws.Cells.Item(11, 15).Value = 0.065  # O11
$ws.Cells.Item(11, 16).Value = 0.255  # P11
$ws.Cells.Item(11, 23).Value = 0.512  # W11
$ws.Cells.Item(11, 24).Value = 0.25  # X11
$ws.Cells.Item(11, 25).Value = 0.5  # Y11
$ws.Cells.Item(11, 35).Value = 0.581  # AI11
$ws.Cells.Item(11, 36).Value = 0.243  # AJ11
$ws.Cells.Item(11, 37).Value = 0.493  # AK11
$ws.Cells.Item(11, 47).Value = 0.558  # AU11
$ws.Cells.Item(11, 48).Value = 0.247  # AV11
$ws.Cells.Item(11, 49).Value = 0.497  # AW11
$ws.Cells.Item(11, 53).Value = 2.186  # BA11
$ws.Cells.Item(11, 54).Value = 0.211  # BB11
$ws.Cells.Item(11, 55).Value = 0.459  # BC11
$ws.Cells.Item(11, 59).Value = 0.674  # BG11
$ws.Cells.Item(11, 60).Value = 0.22  # BH11
$ws.Cells.Item(11, 61).Value = 0.469  # BI11
$ws.Cells.Item(11, 65).Value = 0.814  # BM11
$ws.Cells.Item(11, 66).Value = 0.151  # BN11
$ws.Cells.Item(11, 67).Value = 0.389  # BO11
$ws.Cells.Item(11, 68).Value = 0.729  # BP11
$ws.Cells.Item(11, 69).Value = 0.764  # BQ11
# Row 12
$ws.Cells.Item(12, 5).Value = 1.405  # E12
$ws.Cells.Item(12, 6).Value = 0.836  # F12
$ws.Cells.Item(12, 7).Value = 0.914  # G12
$ws.Cells.Item(12, 14).Value = 1.25  # N12
$ws.Cells.Item(12, 15).Value = 0.287  # O12
$ws.Cells.Item(12, 16).Value = 0.536  # P12
$ws.Cells.Item(12, 23).Value = 1.5  # W12
$ws.Cells.Item(12, 24).Value = 0.432  # X12
$ws.Cells.Item(12, 25).Value = 0.657  # Y12
$ws.Cells.Item(12, 35).Value = 1.6  # AI12
$ws.Cells.Item(12, 36).Value = 1.44  # AJ12
$ws.Cells.Item(12, 37).Value = 1.2  # AK12
$ws.Cells.Item(12, 47).Value = 2.846  # AU12
$ws.Cells.Item(12, 48).Value = 3.361  # AV12
$ws.Cells.Item(12, 49).Value = 1.833  # AW12
$ws.Cells.Item(12, 53).Value = 3.767  # BA12
$ws.Cells.Item(12, 54).Value = 0.44  # BB12
$ws.Cells.Item(12, 55).Value = 0.663  # BC12
$ws.Cells.Item(12, 59).Value = 1.138  # BG12
$ws.Cells.Item(12, 60).Value = 0.188  # BH12
$ws.Cells.Item(12, 61).Value = 0.433  # BI12
$ws.Cells.Item(12, 65).Value = 1.229  # BM12
$ws.Cells.Item(12, 66).Value = 0.233  # BN12
$ws.Cells.Item(12, 67).Value = 0.483  # BO12
$ws.Cells.Item(12, 68).Value = 1.256  # BP12
$ws.Cells.Item(12, 69).Value = 1.241  # BQ12
# Row 13
$ws.Cells.Item(13, 5).Value = 1.405  # E13
$ws.Cells.Item(13, 6).Value = 0.292  # F13
$ws.Cells.Item(13, 7).Value = 0.54  # G13
$ws.Cells.Item(13, 14).Value = 1.737  # N13
$ws.Cells.Item(13, 15).Value = 0.466  # O13
$ws.Cells.Item(13, 16).Value = 0.683  # P13
$ws.Cells.Item(13, 23).Value = 0.985  # W13
$ws.Cells.Item(13, 24).Value = 0.194  # X13
$ws.Cells.Item(13, 25).Value = 0.441  # Y13
$ws.Cells.Item(13, 35).Value = 1.154  # AI13
$ws.Cells.Item(13, 36).Value = 0.303  # AJ13
$ws.Cells.Item(13, 37).Value = 0.551  # AK13
$ws.Cells.Item(13, 47).Value = 2.039  # AU13
$ws.Cells.Item(13, 48).Value = 0.339  # AV13
$ws.Cells.Item(13, 49).Value = 0.582  # AW13
$ws.Cells.Item(13, 53).Value = 2.159  # BA13
$ws.Cells.Item(13, 54).Value = 0.277  # BB13
$ws.Cells.Item(13, 55).Value = 0.527  # BC13
$ws.Cells.Item(13, 59).Value = 0.542  # BG13
$ws.Cells.Item(13, 60).Value = 0.05  # BH13
$ws.Cells.Item(13, 61).Value = 0.224  # BI13
$ws.Cells.Item(13, 65).Value = 0.776  # BM13
$ws.Cells.Item(13, 66).Value = 0.164  # BN13
$ws.Cells.Item(13, 67).Value = 0.404  # BO13
$ws.Cells.Item(13, 68).Value = 0.72  # BP13
$ws.Cells.Item(13, 69).Value = 0.665  # BQ13
